# Update Pais sheet: COVID country-stats refresh (data pulled at 21:13
# instead of 19:56). A few countries' totals overtook their neighbours in
# the "Casos totales" ranking, so those rows' country name shifts along
# with the refreshed numbers:
#   - Sudafrica leap-frogs Pakistan and Italia (rows 14-16)
#   - Congo leap-frogs Sudan del Sur..Lituania (rows 115-121)
#   - Seychelles leap-frogs Barbados (rows 184-185)
# Every other changed row keeps its original country but gets refreshed
# counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "A1"; Value = "Datos actualizados a 10 de Julio de 2020 a las 21:13" },
    @{ Cell = "B4"; Value = 3257696 }, @{ Cell = "C4"; Value = 37697 }, @{ Cell = "D4"; Value = 1438905 }, @{ Cell = "E4"; Value = 1682536 }, @{ Cell = "G4"; Value = 433 }, @{ Cell = "H4"; Value = 136255 },
    @{ Cell = "B6"; Value = 821634 }, @{ Cell = "C6"; Value = 26792 }, @{ Cell = "D6"; Value = 516206 }, @{ Cell = "E6"; Value = 283284 }, @{ Cell = "G6"; Value = 521 }, @{ Cell = "H6"; Value = 22144 },
    @{ Cell = "A14"; Value = "Sudafrica" }, @{ Cell = "B14"; Value = 250687 }, @{ Cell = "C14"; Value = 12348 }, @{ Cell = "D14"; Value = 118232 }, @{ Cell = "E14"; Value = 128595 }, @{ Cell = "G14"; Value = 140 }, @{ Cell = "H14"; Value = 3860 },
    @{ Cell = "A15"; Value = "Pakistan" }, @{ Cell = "B15"; Value = 243599 }, @{ Cell = "C15"; Value = 2751 }, @{ Cell = "D15"; Value = 149092 }, @{ Cell = "E15"; Value = 89449 }, @{ Cell = "G15"; Value = 75 }, @{ Cell = "H15"; Value = 5058 },
    @{ Cell = "A16"; Value = "Italia" }, @{ Cell = "B16"; Value = 242639 }, @{ Cell = "C16"; Value = 276 }, @{ Cell = "D16"; Value = 194273 }, @{ Cell = "E16"; Value = 13428 }, @{ Cell = "G16"; Value = 12 }, @{ Cell = "H16"; Value = 34938 },
    @{ Cell = "B19"; Value = 199363 }, @{ Cell = "C19"; Value = 165 }, @{ Cell = "E19"; Value = 6235 }, @{ Cell = "G19"; Value = 3 }, @{ Cell = "H19"; Value = 9128 },
    @{ Cell = "B21"; Value = 170752 }, @{ Cell = "C21"; Value = 658 }, @{ Cell = "D21"; Value = 78388 }, @{ Cell = "E21"; Value = 62360 }, @{ Cell = "G21"; Value = 25 }, @{ Cell = "H21"; Value = 30004 },
    @{ Cell = "E70"; Value = 4468 }, @{ Cell = "G70"; Value = 3 }, @{ Cell = "H70"; Value = 54 },
    @{ Cell = "B94"; Value = 5203 }, @{ Cell = "C94"; Value = 77 }, @{ Cell = "D94"; Value = 2111 }, @{ Cell = "E94"; Value = 2946 }, @{ Cell = "G94"; Value = 2 }, @{ Cell = "H94"; Value = 146 },
    @{ Cell = "B110"; Value = 2454 }, @{ Cell = "C110"; Value = 300 }, @{ Cell = "E110"; Value = 463 },
    @{ Cell = "A115"; Value = "Congo" }, @{ Cell = "B115"; Value = 2028 }, @{ Cell = "C115"; Value = 207 }, @{ Cell = "D115"; Value = 589 }, @{ Cell = "E115"; Value = 1392 }, @{ Cell = "H115"; Value = 47 },
    @{ Cell = "A116"; Value = "Sudan del Sur" }, @{ Cell = "B116"; Value = 2021 }, @{ Cell = "C116"; Value = 0 }, @{ Cell = "D116"; Value = 333 }, @{ Cell = "E116"; Value = 1650 }, @{ Cell = "H116"; Value = 38 },
    @{ Cell = "A117"; Value = "Estonia" }, @{ Cell = "B117"; Value = 2013 }, @{ Cell = "C117"; Value = 2 }, @{ Cell = "D117"; Value = 1894 }, @{ Cell = "E117"; Value = 50 }, @{ Cell = "H117"; Value = 69 },
    @{ Cell = "A118"; Value = "Zambia" }, @{ Cell = "B118"; Value = 1895 }, @{ Cell = "C118"; Value = 0 }, @{ Cell = "D118"; Value = 1348 }, @{ Cell = "E118"; Value = 505 }, @{ Cell = "H118"; Value = 42 },
    @{ Cell = "A119"; Value = "Islandia" }, @{ Cell = "B119"; Value = 1886 }, @{ Cell = "C119"; Value = 4 }, @{ Cell = "D119"; Value = 1859 }, @{ Cell = "E119"; Value = 17 }, @{ Cell = "H119"; Value = 10 },
    @{ Cell = "A120"; Value = "Eslovaquia" }, @{ Cell = "B120"; Value = 1870 }, @{ Cell = "C120"; Value = 19 }, @{ Cell = "D120"; Value = 1481 }, @{ Cell = "E120"; Value = 361 }, @{ Cell = "H120"; Value = 28 },
    @{ Cell = "A121"; Value = "Lituania" }, @{ Cell = "B121"; Value = 1861 }, @{ Cell = "C121"; Value = 4 }, @{ Cell = "D121"; Value = 1569 }, @{ Cell = "E121"; Value = 213 }, @{ Cell = "H121"; Value = 79 },
    @{ Cell = "B124"; Value = 1613 }, @{ Cell = "C124"; Value = 15 }, @{ Cell = "D124"; Value = 1133 }, @{ Cell = "E124"; Value = 417 },
    @{ Cell = "B146"; Value = 874 }, @{ Cell = "C146"; Value = 1 }, @{ Cell = "D146"; Value = 789 },
    @{ Cell = "B149"; Value = 727 }, @{ Cell = "C149"; Value = 1 }, @{ Cell = "D149"; Value = 284 },
    @{ Cell = "B155"; Value = 668 }, @{ Cell = "C155"; Value = 53 }, @{ Cell = "E155"; Value = 642 },
    @{ Cell = "A184"; Value = "Seychelles" }, @{ Cell = "B184"; Value = 100 }, @{ Cell = "C184"; Value = 6 }, @{ Cell = "D184"; Value = 11 }, @{ Cell = "E184"; Value = 89 }, @{ Cell = "H184"; Value = 0 },
    @{ Cell = "A185"; Value = "Barbados" }, @{ Cell = "B185"; Value = 98 }, @{ Cell = "D185"; Value = 90 }, @{ Cell = "E185"; Value = 1 }, @{ Cell = "H185"; Value = 7 }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}
